# Revert "adding term 2.0.0"
#
# 1. Restore the previous metadata values (Version, Date, Contact) on the
#    "Metadata" sheet.
# 2. Restore the previous "Value" code (descendent-of = "D") on the
#    "Include from FSIII" sheet.
# 3. Remove the duplicated "Include from FSIII 2" worksheet that the
#    reverted commit had added.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

$incl = $wb.Worksheets.Item("Include from FSIII")
$incl.Range("C2").Value = "D"

$dup = $wb.Worksheets.Item("Include from FSIII 2")
$dup.Delete()
